$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range('D2').Value = '47.196.35'
$ws.Range('D3').Value = '2.484.33'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.83'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.04'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.08'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.39'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.16'
$ws.Range('D15').Value = '2.873.67'
$ws.Range('D16').Value = '2.488.57'
$ws.Range('D18').Value = '47.121.39'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.47'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.62'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.76'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.46'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '245.83'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.64'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.29'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.94'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.55'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.77'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.27'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.76'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.93'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.03'
$ws.Range('D45').Value = '1.994.61'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.80'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.09'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.09'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.59'

# --- Column E (Volume 1h) updates ---
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  +3.53%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('E19').Value = '  +6.32%  '
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('E22').Value = '  +15.45%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('E28').Value = '  +3.90%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +2.55%  '
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('E40').Value = '  +8.16%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('E47').Value = '  -4.55%  '
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('E49').Value = '  -1.23%  '
$ws.Range('E50').Value = '  -4.62%  '
$ws.Range('E51').Value = '  +3.31%  '
